$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EPIDEMI PENYAKIT")
$ws.Activate()

$codes = @(
    "53.06.13.2021",
    "53.06.13.2020",
    "53.06.13.2019",
    "53.06.13.2018",
    "53.06.13.2017",
    "53.06.13.2016",
    "53.06.13.2015",
    "53.06.13.2014"
)

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $codes[$i]
}

$ws.Range("F9").Select()
